$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1682.1428
$ws.Range("J17").Value = 1682.1428
$ws.Range("L17").Value = 5046.428400000001
$ws.Range("N17").Value = -5382.428400000001
$ws.Range("H48").Value = 3860.75
$ws.Range("J48").Value = 4722
$ws.Range("L48").Value = 14166
$ws.Range("N48").Value = -14750
$ws.Range("H56").Value = 3860.75
$ws.Range("J56").Value = 4722
$ws.Range("L56").Value = 14166
$ws.Range("N56").Value = -15234
$ws.Range("H137").Value = 6602.95
$ws.Range("J137").Value = 7729.1665
$ws.Range("L137").Value = 23187.4995
$ws.Range("N137").Value = -28287.4995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2809.652
$ws.Range("I94").Value = 1194
$ws.Range("K94").Value = 1194
$ws.Range("M94").Value = -743
$ws.Range("N135").Value = -58139.5
$ws.Range("H135").Value = 47999.5
$ws.Range("J135").Value = 47999.5
$ws.Range("L135").Value = 47999.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2311.7222
$ws.Range("I31").Value = 1280.8
$ws.Range("K31").Value = 1280.8
$ws.Range("M31").Value = -985.8
$ws.Range("H34").Value = 2311.7222
$ws.Range("I34").Value = 1280.8
$ws.Range("K34").Value = 1280.8
$ws.Range("M34").Value = -1078.8
$ws.Range("H58").Value = 1543.0968
$ws.Range("I58").Value = 1459.091
$ws.Range("J58").Value = 1748.4445
$ws.Range("K58").Value = 1459.091
$ws.Range("L58").Value = 1748.4445
$ws.Range("M58").Value = -1256.091
$ws.Range("N58").Value = -2154.4445
$ws.Range("H62").Value = 22042.572
$ws.Range("I62").Value = 2749.5
$ws.Range("J62").Value = 29759.8
$ws.Range("K62").Value = 2749.5
$ws.Range("L62").Value = 29759.8
$ws.Range("M62").Value = -2125.5
$ws.Range("N62").Value = -31007.8
$ws.Range("H65").Value = 22042.572
$ws.Range("I65").Value = 2749.5
$ws.Range("J65").Value = 29759.8
$ws.Range("K65").Value = 13747.5
$ws.Range("L65").Value = 148799
$ws.Range("M65").Value = -10627.5
$ws.Range("N65").Value = -155039
$ws.Range("H99").Value = 2164.625
$ws.Range("I99").Value = 2009.4
$ws.Range("J99").Value = 2423.3333
$ws.Range("K99").Value = 2009.4
$ws.Range("L99").Value = 2423.3333
$ws.Range("M99").Value = -511.4000000000001
$ws.Range("N99").Value = -5419.3333
$ws.Range("H105").Value = 2786.25
$ws.Range("I105").Value = 2170.6667
$ws.Range("J105").Value = 4633
$ws.Range("K105").Value = 2170.6667
$ws.Range("L105").Value = 4633
$ws.Range("M105").Value = -423.6667000000002
$ws.Range("N105").Value = -8127
$ws.Range("H122").Value = 2264.3333
$ws.Range("I122").Value = 2149.5
$ws.Range("J122").Value = 2494
$ws.Range("K122").Value = 6448.5
$ws.Range("L122").Value = 7482
$ws.Range("M122").Value = -3998.5
$ws.Range("N122").Value = -12382
$ws.Range("H126").Value = 2164.625
$ws.Range("I126").Value = 2009.4
$ws.Range("J126").Value = 2423.3333
$ws.Range("K126").Value = 6028.200000000001
$ws.Range("L126").Value = 7269.999899999999
$ws.Range("M126").Value = -3558.200000000001
$ws.Range("N126").Value = -12209.9999
$ws.Range("H132").Value = 4562.087
$ws.Range("I132").Value = 4477.8887
$ws.Range("K132").Value = 13433.6661
$ws.Range("M132").Value = -10903.6661
$ws.Range("H134").Value = 2328634
$ws.Range("I134").Value = 2351.139
$ws.Range("K134").Value = 7053.417
$ws.Range("M134").Value = -4518.417
$ws.Range("H136").Value = 1543.0968
$ws.Range("I136").Value = 1459.091
$ws.Range("J136").Value = 1748.4445
$ws.Range("K136").Value = 4377.272999999999
$ws.Range("L136").Value = 5245.333500000001
$ws.Range("M136").Value = -1827.272999999999
$ws.Range("N136").Value = -10345.3335

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M12").ClearContents()
$ws.Range("H12").Value = 1696.6
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1696.6
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 5089.799999999999
$ws.Range("N12").Value = -5435.799999999999
$ws.Range("H40").Value = 150.73914
$ws.Range("I40").Value = 162.5
$ws.Range("K40").Value = 650
$ws.Range("M40").Value = -581
$ws.Range("M97").ClearContents()
$ws.Range("H97").Value = 50279.5
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("H103").Value = 150.61539
$ws.Range("J103").Value = 112.25
$ws.Range("L103").Value = 336.75
$ws.Range("N103").Value = -2094.75
$ws.Range("H117").Value = 2826.7144
$ws.Range("J117").Value = 4353.375
$ws.Range("L117").Value = 13060.125
$ws.Range("N117").Value = -19944.125
$ws.Range("H130").Value = 5032.0835
$ws.Range("J130").Value = 5032.909
$ws.Range("L130").Value = 15098.727
$ws.Range("N130").Value = -25138.727
$ws.Range("H134").Value = 6054.1113
$ws.Range("I134").Value = 5625.2856
$ws.Range("K134").Value = 16875.8568
$ws.Range("M134").Value = -11805.8568

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2975.4644
$ws.Range("I122").Value = 3022.7727
$ws.Range("K122").Value = 9068.3181
$ws.Range("M122").Value = -6618.3181
$ws.Range("H132").Value = 1162.375
$ws.Range("I132").Value = 1042.7142
$ws.Range("K132").Value = 3128.1426
$ws.Range("M132").Value = -598.1425999999997

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 3305.5
$ws.Range("I53").Value = 3079
$ws.Range("J53").Value = 3683
$ws.Range("K53").Value = 3079
$ws.Range("L53").Value = 3683
$ws.Range("M53").Value = -2561
$ws.Range("N53").Value = -4719
$ws.Range("N61").ClearContents()
$ws.Range("H61").Value = 2334.5715
$ws.Range("I61").Value = 2334.5715
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2334.5715
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2132.5715
$ws.Range("N113").ClearContents()
$ws.Range("H113").Value = 2334.5715
$ws.Range("I113").Value = 2334.5715
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2334.5715
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -164.5715
$ws.Range("H122").Value = 3402.2144
$ws.Range("J122").Value = 3529.4443
$ws.Range("L122").Value = 10588.3329
$ws.Range("N122").Value = -15488.3329

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 70123.30499999999
$ws.Range("J41").Value = 76986.91
$ws.Range("L41").Value = 76986.91
$ws.Range("N41").Value = -77766.91
$ws.Range("H70").Value = 31087.666
$ws.Range("J70").Value = 31087.666
$ws.Range("L70").Value = 31087.666
$ws.Range("N70").Value = -31717.666
$ws.Range("H73").Value = 31087.666
$ws.Range("J73").Value = 31087.666
$ws.Range("L73").Value = 31087.666
$ws.Range("N73").Value = -33271.666
$ws.Range("H133").Value = 48000
$ws.Range("J133").Value = 48000
$ws.Range("L133").Value = 48000
$ws.Range("N133").Value = -58120
